# Auto update: 2025-12-05 13:19:50
# Refresh the hedging/gold analysis sheet: rows 2 and 3 swap identity
# (Gold Feb 26 / GC=F <-> StreetTRACKS Gold Shares / GLD) and every row's
# market-derived numbers (close, RSI, 5d return, probabilities, final
# score, macro score) are refreshed with the latest pulled values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> now StreetTRACKS Gold Shares / GLD
$ws.Range("B2").Value = "StreetTRACKS Gold Shares"
$ws.Range("C2").Value = "GLD"
$ws.Range("D2").Value = 387.13
$ws.Range("E2").Value = 56.3
$ws.Range("F2").Value = 1.05
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 73
$ws.Range("I2").Value = 83
$ws.Range("J2").Value = 96
$ws.Range("K2").Value = 67.7
$ws.Range("N2").Value = 54.83846622768671

# Row 3 -> now Gold Feb 26 / GC=F
$ws.Range("B3").Value = "Gold Feb 26"
$ws.Range("C3").Value = "GC=F"
$ws.Range("D3").Value = 4242
$ws.Range("E3").Value = 56.1
$ws.Range("F3").Value = 1.84
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 66
$ws.Range("I3").Value = 80
$ws.Range("J3").Value = 83
$ws.Range("K3").Value = 66.5
$ws.Range("N3").Value = 54.83846622768671

# Row 4 -> Newmont Corporation / NEM (identity unchanged, scores refreshed)
$ws.Range("K4").Value = 66.5
$ws.Range("N4").Value = 54.83846622768671
